# Update employment alignment targets
# ---------------------------------------------------------------------------
# Updates the B2:B19 alignment-target series on the three "employment_*"
# sheets, moves the active-cell/selection on several sheets, and switches
# the active tab from "disability" to "employment_couples".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. employment_smales: new B2:B19 values + new selection (B14:B19)
# ---------------------------------------------------------------------------
$wsSM = $wb.Worksheets.Item("employment_smales")

$smalesValues = @(
    0.63684130000000005,
    0.63684130000000005,
    0.63086129999999996,
    0.62820810000000005,
    0.63444250000000002,
    0.64983950000000001,
    0.66639119999999996,
    0.67221770000000003,
    0.67153039999999997,
    0.67241660000000003,
    0.69137210000000004,
    0.703434,
    0.69511619999999996,
    0.69436249999999999,
    0.69360880000000003,
    0.69285509999999995,
    0.69210139999999998,
    0.69134770000000001
)

for ($i = 0; $i -lt $smalesValues.Length; $i++) {
    $row = 2 + $i
    $wsSM.Cells.Item($row, 2).Value = $smalesValues[$i]
}

# ---------------------------------------------------------------------------
# 2. employment_sfemales: new B2:B19 values + new selection (B14:B19)
# ---------------------------------------------------------------------------
$wsSF = $wb.Worksheets.Item("employment_sfemales")

$sfemalesValues = @(
    0.39854529999999999,
    0.39854529999999999,
    0.4029663,
    0.39655180000000001,
    0.4070319,
    0.40790100000000001,
    0.41726039999999998,
    0.42810779999999998,
    0.42215740000000002,
    0.41901169999999999,
    0.42423159999999999,
    0.41884729999999998,
    0.42949969999999998,
    0.42430040000000002,
    0.4191011,
    0.41390179999999999,
    0.40870250000000002,
    0.40350320000000001
)

for ($i = 0; $i -lt $sfemalesValues.Length; $i++) {
    $row = 2 + $i
    $wsSF.Cells.Item($row, 2).Value = $sfemalesValues[$i]
}

# ---------------------------------------------------------------------------
# 3. employment_couples: new B2:B19 values + new selection (A2), becomes
#    the active sheet/tab (previously "disability" was active).
# ---------------------------------------------------------------------------
$wsEC = $wb.Worksheets.Item("employment_couples")

$couplesValues = @(
    0.7511234,
    0.7511234,
    0.74668140000000005,
    0.74351520000000004,
    0.74422900000000003,
    0.73707429999999996,
    0.74269050000000003,
    0.74164249999999998,
    0.73936999999999997,
    0.7404522,
    0.72219089999999997,
    0.72516029999999998,
    0.74010039999999999,
    0.73162090000000002,
    0.738754466666667,
    0.74198476666666702,
    0.74521506666666704,
    0.74844536666666694
)

for ($i = 0; $i -lt $couplesValues.Length; $i++) {
    $row = 2 + $i
    $wsEC.Cells.Item($row, 2).Value = $couplesValues[$i]
}

# ---------------------------------------------------------------------------
# 4. Selections: restore/update the per-sheet selection rectangles without
#    disturbing which sheet ends up active (that is set last, below).
# ---------------------------------------------------------------------------
$wsSM.Activate() | Out-Null
$wsSM.Range("B14:B19").Select() | Out-Null

$wsSF.Activate() | Out-Null
$wsSF.Range("B14:B19").Select() | Out-Null

# disability: drop the old tabSelected flag by moving focus away from it,
# leave its own selection (B1) untouched.
$wsDis = $wb.Worksheets.Item("disability")
$wsDis.Activate() | Out-Null
$wsDis.Range("B1").Select() | Out-Null

# employment_couples becomes the active/selected tab, with A2 selected.
$wsEC.Activate() | Out-Null
$wsEC.Range("A2").Select() | Out-Null
